$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1252.4789
$ws.Range("I15").Value = 1252.4789
$ws.Range("K15").Value = 3757.4367
$ws.Range("M15").Value = -3588.4367
$ws.Range("H33").Value = 104.875
$ws.Range("I33").Value = 80
$ws.Range("J33").Value = 113.166664
$ws.Range("K33").Value = 80
$ws.Range("L33").Value = 113.166664
$ws.Range("M33").Value = 149
$ws.Range("N33").Value = -571.166664
$ws.Range("H38").Value = 440.6
$ws.Range("I38").Value = 440.6
$ws.Range("K38").Value = 1321.8
$ws.Range("M38").Value = -949.8000000000002
$ws.Range("H42").Value = 595.1111
$ws.Range("I42").Value = 94
$ws.Range("J42").Value = 738.2857
$ws.Range("K42").Value = 282
$ws.Range("L42").Value = 2214.8571
$ws.Range("M42").Value = -52
$ws.Range("N42").Value = -2674.8571
$ws.Range("H87").Value = 57998.332
$ws.Range("J87").Value = 57998.332
$ws.Range("L87").Value = 57998.332
$ws.Range("N87").Value = -60494.332
$ws.Range("H90").Value = 57998.332
$ws.Range("J90").Value = 57998.332
$ws.Range("L90").Value = 173994.996
$ws.Range("N90").Value = -186474.996
$ws.Range("H106").Value = 7669.6665
$ws.Range("I106").Value = 3997
$ws.Range("J106").Value = 9506
$ws.Range("K106").Value = 3997
$ws.Range("L106").Value = 9506
$ws.Range("M106").Value = -3366
$ws.Range("N106").Value = -10768
$ws.Range("H113").Value = 24789.846
$ws.Range("I113").Value = 38683.5
$ws.Range("J113").Value = 2560
$ws.Range("K113").Value = 38683.5
$ws.Range("L113").Value = 2560
$ws.Range("M113").Value = -35429.5
$ws.Range("N113").Value = -9068
$ws.Range("H137").Value = 31667.969
$ws.Range("J137").Value = 73008
$ws.Range("L137").Value = 219024
$ws.Range("N137").Value = -224124
$ws.Range("H139").Value = 60641.6
$ws.Range("J139").Value = 60641.6
$ws.Range("L139").Value = 60641.6
$ws.Range("N139").Value = -70921.60000000001
$ws.Range("H140").Value = 83578.375
$ws.Range("J140").Value = 83578.375
$ws.Range("L140").Value = 83578.375
$ws.Range("N140").Value = -93938.375
$ws.Range("H141").Value = 801809.7
$ws.Range("I141").Value = 875902.6
$ws.Range("J141").Value = 11485.333
$ws.Range("K141").Value = 2627707.8
$ws.Range("L141").Value = 34455.999
$ws.Range("M141").Value = -2622527.8
$ws.Range("N141").Value = -44815.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3867.695
$ws.Range("I32").Value = 3328.5972
$ws.Range("K32").Value = 3328.5972
$ws.Range("M32").Value = -3041.5972

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23651.4
$ws.Range("I82").Value = 17564.25
$ws.Range("J82").Value = 48000
$ws.Range("K82").Value = 17564.25
$ws.Range("L82").Value = 48000
$ws.Range("M82").Value = -17181.25
$ws.Range("N82").Value = -48766
$ws.Range("H85").Value = 23651.4
$ws.Range("I85").Value = 17564.25
$ws.Range("J85").Value = 48000
$ws.Range("K85").Value = 17564.25
$ws.Range("L85").Value = 48000
$ws.Range("M85").Value = -16238.25
$ws.Range("N85").Value = -50652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2001.64
$ws.Range("I31").Value = 1462.25
$ws.Range("K31").Value = 1462.25
$ws.Range("M31").Value = -1167.25
$ws.Range("H34").Value = 2001.64
$ws.Range("I34").Value = 1462.25
$ws.Range("K34").Value = 1462.25
$ws.Range("M34").Value = -1260.25
$ws.Range("H122").Value = 4692.077
$ws.Range("I122").Value = 3398.5
$ws.Range("J122").Value = 9004
$ws.Range("K122").Value = 10195.5
$ws.Range("L122").Value = 27012
$ws.Range("M122").Value = -7745.5
$ws.Range("N122").Value = -31912
$ws.Range("H132").Value = 2263.4546
$ws.Range("I132").Value = 1523.0834
$ws.Range("J132").Value = 4237.778
$ws.Range("K132").Value = 4569.2502
$ws.Range("L132").Value = 12713.334
$ws.Range("M132").Value = -2039.2502
$ws.Range("N132").Value = -17773.334
$ws.Range("H134").Value = 1349.6338
$ws.Range("I134").Value = 1224.0518
$ws.Range("J134").Value = 1909.9231
$ws.Range("K134").Value = 3672.1554
$ws.Range("L134").Value = 5729.7693
$ws.Range("M134").Value = -1137.1554
$ws.Range("N134").Value = -10799.7693
$ws.Range("H138").Value = 98136.37
$ws.Range("J138").Value = 98136.37
$ws.Range("L138").Value = 98136.37
$ws.Range("N138").Value = -108416.37

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1597.6
$ws.Range("I80").Value = 1349
$ws.Range("J80").Value = 1763.3334
$ws.Range("K80").Value = 4047
$ws.Range("L80").Value = 5290.0002
$ws.Range("M80").Value = -3111
$ws.Range("N80").Value = -7162.0002
$ws.Range("H83").Value = 1597.6
$ws.Range("I83").Value = 1349
$ws.Range("J83").Value = 1763.3334
$ws.Range("K83").Value = 12141
$ws.Range("L83").Value = 15870.0006
$ws.Range("M83").Value = -7461
$ws.Range("N83").Value = -25230.0006
$ws.Range("H131").Value = 15179398
$ws.Range("I131").Value = 45455056
$ws.Range("J131").Value = 41568.137
$ws.Range("K131").Value = 136365168
$ws.Range("L131").Value = 124704.411
$ws.Range("M131").Value = -136360128
$ws.Range("N131").Value = -134784.411
$ws.Range("H132").Value = 1286.55
$ws.Range("J132").Value = 1503.0714
$ws.Range("L132").Value = 13527.6426
$ws.Range("N132").Value = -18587.6426
$ws.Range("H133").Value = 17861054
$ws.Range("I133").Value = 41669130
$ws.Range("K133").Value = 125007390
$ws.Range("M133").Value = -125002330
$ws.Range("H134").Value = 25450.883
$ws.Range("I134").Value = 33657.676
$ws.Range("J134").Value = 4250
$ws.Range("K134").Value = 100973.028
$ws.Range("L134").Value = 12750
$ws.Range("M134").Value = -95903.02799999999
$ws.Range("N134").Value = -22890
$ws.Range("H136").Value = 1831.3572
$ws.Range("I136").Value = 1363.9
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4091.7
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 1008.3
$ws.Range("N136").Value = -19200
$ws.Range("H137").Value = 3495.3333
$ws.Range("I137").Value = 1799.6364
$ws.Range("J137").Value = 4930.154
$ws.Range("K137").Value = 5398.9092
$ws.Range("L137").Value = 14790.462
$ws.Range("M137").Value = -298.9092000000001
$ws.Range("N137").Value = -24990.462
$ws.Range("H138").Value = 2511.6428
$ws.Range("I138").Value = 2219.1667
$ws.Range("J138").Value = 4266.5
$ws.Range("K138").Value = 6657.500100000001
$ws.Range("L138").Value = 12799.5
$ws.Range("M138").Value = -1517.500100000001
$ws.Range("N138").Value = -23079.5
$ws.Range("H139").Value = 5288.923
$ws.Range("I139").Value = 5420.56
$ws.Range("J139").Value = 1998
$ws.Range("K139").Value = 16261.68
$ws.Range("L139").Value = 5994
$ws.Range("M139").Value = -11121.68
$ws.Range("N139").Value = -16274
$ws.Range("H140").Value = 2226.3333
$ws.Range("I140").Value = 1193.1177
$ws.Range("J140").Value = 2742.9412
$ws.Range("K140").Value = 3579.3531
$ws.Range("L140").Value = 8228.8236
$ws.Range("M140").Value = 1600.6469
$ws.Range("N140").Value = -18588.8236
$ws.Range("H141").Value = 2878.6924
$ws.Range("I141").Value = 2910.25
$ws.Range("K141").Value = 8730.75
$ws.Range("M141").Value = -3550.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1744150.6
$ws.Range("I14").Value = 2241609.8
$ws.Range("J14").Value = 500502.5
$ws.Range("K14").Value = 2241609.8
$ws.Range("L14").Value = 500502.5
$ws.Range("M14").Value = -2241441.8
$ws.Range("N14").Value = -500838.5
$ws.Range("H43").Value = 6347
$ws.Range("I43").Value = 2129.3333
$ws.Range("K43").Value = 2129.3333
$ws.Range("M43").Value = -1978.3333
$ws.Range("H46").Value = 21533.334
$ws.Range("J46").Value = 21533.334
$ws.Range("L46").Value = 21533.334
$ws.Range("N46").Value = -21845.334
$ws.Range("H80").Value = 700
$ws.Range("I80").Value = 700
$ws.Range("K80").Value = 700
$ws.Range("M80").Value = 298
$ws.Range("H83").Value = 700
$ws.Range("I83").Value = 700
$ws.Range("K83").Value = 3500
$ws.Range("M83").Value = 1492
$ws.Range("H127").Value = 36639
$ws.Range("J127").Value = 36639
$ws.Range("L127").Value = 36639
$ws.Range("N127").Value = -46559
$ws.Range("H132").Value = 1101528.6
$ws.Range("I132").Value = 1833633.2
$ws.Range("J132").Value = 3371.7144
$ws.Range("K132").Value = 5500899.6
$ws.Range("L132").Value = 10115.1432
$ws.Range("M132").Value = -5498369.6
$ws.Range("N132").Value = -15175.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3872.3845
$ws.Range("I7").Value = 3104.5557
$ws.Range("J7").Value = 5600
$ws.Range("K7").Value = 3104.5557
$ws.Range("L7").Value = 5600
$ws.Range("M7").Value = -2992.5557
$ws.Range("N7").Value = -5824
$ws.Range("H82").Value = 2002.1
$ws.Range("I82").Value = 1475.1666
$ws.Range("J82").Value = 2792.5
$ws.Range("K82").Value = 1475.1666
$ws.Range("L82").Value = 2792.5
$ws.Range("M82").Value = -1114.1666
$ws.Range("N82").Value = -3514.5
$ws.Range("H85").Value = 2002.1
$ws.Range("I85").Value = 1475.1666
$ws.Range("J85").Value = 2792.5
$ws.Range("K85").Value = 1475.1666
$ws.Range("L85").Value = 2792.5
$ws.Range("M85").Value = -227.1666
$ws.Range("N85").Value = -5288.5
$ws.Range("H126").Value = 3872.3845
$ws.Range("I126").Value = 3104.5557
$ws.Range("J126").Value = 5600
$ws.Range("K126").Value = 9313.667099999999
$ws.Range("L126").Value = 16800
$ws.Range("M126").Value = -6843.667099999999
$ws.Range("N126").Value = -21740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 14853.4
$ws.Range("J74").Value = 14853.4
$ws.Range("L74").Value = 14853.4
$ws.Range("N74").Value = -16725.4
$ws.Range("H77").Value = 14853.4
$ws.Range("J77").Value = 14853.4
$ws.Range("L77").Value = 44560.2
$ws.Range("N77").Value = -53920.2
$ws.Range("H96").Value = 10786.667
$ws.Range("I96").Value = 3345.25
$ws.Range("J96").Value = 16739.8
$ws.Range("K96").Value = 3345.25
$ws.Range("L96").Value = 16739.8
$ws.Range("M96").Value = -1972.25
$ws.Range("N96").Value = -19485.8
$ws.Range("H107").Value = 999.4167
$ws.Range("I107").Value = 550.2
$ws.Range("K107").Value = 1650.6
$ws.Range("M107").Value = 269.3999999999999
$ws.Range("H132").Value = 1905.6216
$ws.Range("I132").Value = 1451.8334
$ws.Range("J132").Value = 2743.3845
$ws.Range("K132").Value = 4355.5002
$ws.Range("L132").Value = 8230.1535
$ws.Range("M132").Value = -1825.5002
$ws.Range("N132").Value = -13290.1535
